$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
  2 = @(44315, 100, 20000, 21000, 20500, 1025)
  3 = @(44410, 200, 20000, 21000, 20500, 1025)
  4 = @(44417, 160, 20000, 21000, 20500, 1025)
  5 = @(44365, 100, 20000, 21000, 20500, 1025)
  6 = @(44333, 100, 19500, 20000, 19750, 988)
  7 = @(44462, 100, 19500, 20000, 19750, 988)
  8 = @(44364, 140, 20000, 21000, 20500, 1025)
  9 = @(44467, 200, 20000, 21000, 20500, 1025)
  10 = @(44301, 100, 18000, 19000, 18500, 925)
  11 = @(44445, 160, 20000, 21000, 20500, 1025)
  12 = @(44420, 160, 20000, 21000, 20500, 1025)
  13 = @(44428, 100, 20000, 21000, 20500, 1025)
  14 = @(44431, 160, 21000, 22000, 21500, 1075)
  15 = @(44466, 100, 20000, 21000, 20500, 1025)
  16 = @(44427, 200, 20000, 21000, 20500, 1025)
  17 = @(44434, 100, 20000, 21000, 20500, 1025)
  18 = @(44448, 100, 20000, 21000, 20500, 1025)
  19 = @(44441, 160, 20000, 21000, 20500, 1025)
  20 = @(44343, 100, 19500, 20000, 19750, 988)
  21 = @(44407, 160, 20000, 21000, 20500, 1025)
  22 = @(44350, 160, 19000, 20000, 19500, 975)
  23 = @(44435, 260, 20000, 22000, 21115, 1056)
  24 = @(44335, 200, 19000, 20000, 19500, 975)
  25 = @(44474, 200, 19000, 20000, 19500, 975)
  26 = @(44418, 200, 20000, 21000, 20500, 1025)
  27 = @(44473, 40, 19500, 20000, 19750, 988)
  28 = @(44326, 160, 19500, 20000, 19750, 988)
  29 = @(44336, 100, 19500, 20000, 19750, 988)
  30 = @(44442, 140, 20000, 21000, 20500, 1025)
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  $ws.Range("D" + $r).Value = $vals[0]
  $ws.Range("M" + $r).Value = $vals[1]
  $ws.Range("N" + $r).Value = $vals[2]
  $ws.Range("O" + $r).Value = $vals[3]
  $ws.Range("P" + $r).Value = $vals[4]
  $ws.Range("S" + $r).Value = $vals[5]
}

Write-Host "done"